$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.489.62"
$ws.Range("E2").Value = "  +0.81%  "

# Row 3
$ws.Range("D3").Value = "1.925.16"
$ws.Range("E3").Value = "  +1.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("E5").Value = "  +12.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "254.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.76%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0754"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.63%  "

# Row 12
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("D13").Value = "2.205.92"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.00%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.723"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.62%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.934.95"

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.63%  "

# Row 18
$ws.Range("D18").Value = "35.502.83"
$ws.Range("E18").Value = "  +0.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.92%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0840"
$ws.Range("E20").Value = "  +3.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.83%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.62%  "

# Row 24
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.13%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.79%  "

# Row 29
$ws.Range("E29").Value = "  +6.88%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.11%  "

# Row 31
$ws.Range("D31").Value = "4.128.12"
$ws.Range("E31").Value = "  +19.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.95%  "

# Row 33
$ws.Range("B33").Value = "TrustWalletToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +25.98%  "

# Row 34
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0583"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.85%  "

# Row 37
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.918"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.47%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.17%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0211"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0657"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.25%  "

# Row 45
$ws.Range("E45").Value = "  +5.34%  "

# Row 46
$ws.Range("D46").Value = "1.350.47"
$ws.Range("E46").Value = "  +1.16%  "

# Row 47
$ws.Range("E47").Value = "  +0.88%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.40%  "

# Row 49
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.81%  "
